# Updates the "cryptos" price/volume table with the latest scraped values.
# Price values in column D are plain text (they use '.' both as thousands
# separator and decimal separator, e.g. "23.546.56"), so values that would
# otherwise be auto-parsed as a number by Excel are written with a leading
# apostrophe to force a text quote-prefix and keep the exact original text
# (including any trailing zeros) instead of being normalized as a float.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.546.56"
$ws.Range("E2").Value = "  +1.24%  "
$ws.Range("D3").Value = "1.655.75"
$ws.Range("E3").Value = "  +2.49%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("D6").Value = "'302.48"
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("D7").Value = "'0.3837"
$ws.Range("E7").Value = "  +1.50%  "
$ws.Range("D8").Value = "'51.21"
$ws.Range("E8").Value = "  -0.82%  "
$ws.Range("D9").Value = "'0.3593"
$ws.Range("E9").Value = "  +1.74%  "
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "'1.241"
$ws.Range("E10").Value = "  +2.85%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "'0.08190"
$ws.Range("E11").Value = "  +0.96%  "
$ws.Range("D12").Value = "'1.002"
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").Value = "'22.41"
$ws.Range("E13").Value = "  +0.88%  "
$ws.Range("D14").Value = "'6.489"
$ws.Range("E14").Value = "  +1.97%  "
$ws.Range("D15").Value = "'7.494"
$ws.Range("E15").Value = "  +3.02%  "
$ws.Range("D16").Value = "'0.00001220"
$ws.Range("E16").Value = "  +0.65%  "
$ws.Range("D17").Value = "1.653.50"
$ws.Range("E17").Value = "  +3.63%  "
$ws.Range("E18").Value = "  +3.54%  "
$ws.Range("D19").Value = "'0.06972"
$ws.Range("E19").Value = "  +0.87%  "
$ws.Range("D20").Value = "'6.814"
$ws.Range("E20").Value = "  +5.18%  "
$ws.Range("E21").Value = "  +2.50%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("E23").Value = "  +2.71%  "
$ws.Range("D24").Value = "23.549.55"
$ws.Range("E24").Value = "  +1.28%  "
$ws.Range("D25").Value = "'2.509"
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("D26").Value = "'3.005"
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("E27").Value = "  +1.59%  "
$ws.Range("D28").Value = "'152.06"
$ws.Range("E28").Value = "  +0.60%  "
$ws.Range("D29").Value = "'5.241"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").Value = "'133.73"
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'7.210"
$ws.Range("E31").Value = "  +11.41%  "
$ws.Range("B32").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C32").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D32").Value = "1.833.81"
$ws.Range("E32").Value = "  +3.09%  "
$ws.Range("D33").Value = "'2.251"
$ws.Range("E33").Value = "  +7.46%  "
$ws.Range("D34").Value = "'12.19"
$ws.Range("E34").Value = "  +7.48%  "
$ws.Range("D35").Value = "'1.057"
$ws.Range("E35").Value = "  -0.89%  "
$ws.Range("D36").Value = "'0.02805"
$ws.Range("E36").Value = "  +3.57%  "
$ws.Range("D37").Value = "'6.125"
$ws.Range("E37").Value = "  +4.71%  "
$ws.Range("D38").Value = "'0.2496"
$ws.Range("E38").Value = "  +1.70%  "
$ws.Range("D39").Value = "'0.08778"
$ws.Range("E39").Value = "  +0.94%  "
$ws.Range("D40").Value = "'0.07017"
$ws.Range("E40").Value = "  +1.11%  "
$ws.Range("D41").Value = "'13.23"
$ws.Range("E41").Value = "  +10.70%  "
$ws.Range("D42").Value = "'0.7016"
$ws.Range("E42").Value = "  +1.91%  "
$ws.Range("D43").Value = "'1.334"
$ws.Range("E43").Value = "  +0.80%  "
$ws.Range("D44").Value = "'15.95"
$ws.Range("E44").Value = "  +4.30%  "
$ws.Range("D45").Value = "'0.6529"
$ws.Range("E45").Value = "  +3.35%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").Value = "'2.306"
$ws.Range("E47").Value = "  +2.41%  "
$ws.Range("D48").Value = "'3.956"
$ws.Range("E48").Value = "  +0.39%  "
$ws.Range("D49").Value = "'0.07913"
$ws.Range("E49").Value = "  +0.55%  "
$ws.Range("D50").Value = "'127.89"
$ws.Range("E50").Value = "  +0.43%  "
$ws.Range("E51").Value = "  +1.87%  "
